$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets: Arkusz1 -> EU, Arkusz2 -> US
$ws1.Name = "EU"
$ws2.Name = "US"

# --- EU sheet (was Arkusz1): drop row 4, add "time zone" column A, and
#     replace the D/E date-label columns with plain numbers ---
$ws1.Rows.Item(4).Delete()

$ws1.Range("A1").Value = "fr/fr/"
$ws1.Range("A2").Value = "pl/pl/"
$ws1.Range("A3").Value = "es/es/"

$ws1.Range("D1").Value = 20
$ws1.Range("E1").Value = 30
$ws1.Range("D2").Value = 11
$ws1.Range("E2").Value = 22
$ws1.Range("D3").Value = 11
$ws1.Range("E3").Value = 22

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- US sheet (was Arkusz2): fill out full B:E columns + a third row ---
$ws2.Range("A1").Value = "us/en/"
$ws2.Range("B1").Value = "WAW"
$ws2.Range("C1").Value = "JFK"
$ws2.Range("D1").Value = 13
$ws2.Range("E1").Value = 23

$ws2.Range("A2").Value = "us/en/"
$ws2.Range("B2").Value = "WAW"
$ws2.Range("C2").Value = "JFK"
$ws2.Range("D2").Value = 11
$ws2.Range("E2").Value = 22

$ws2.Range("A3").Value = "us/en/"
$ws2.Range("B3").Value = "WAW"
$ws2.Range("C3").Value = "JFK"
$ws2.Range("D3").Value = 11
$ws2.Range("E3").Value = 22

# --- Selections / active tab: EU becomes the active/selected sheet ---
$ws2.Range("A1:E3").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("H8").Select() | Out-Null
